# Update existing rows 2-13 and add new rows 14-16 for the weekly data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 7
$ws.Cells.Item(2,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2,3).Value = "Ñuble"
$ws.Cells.Item(2,4).Value = 44698
$ws.Cells.Item(2,5).Value = 16
$ws.Cells.Item(2,6).Value = "Fruta"
$ws.Cells.Item(2,7).Value = 100107
$ws.Cells.Item(2,8).Value = "Otros"
$ws.Cells.Item(2,9).Value = 100107001
$ws.Cells.Item(2,10).Value = "Caqui"
$ws.Cells.Item(2,11).Value = "Mankaki"
$ws.Cells.Item(2,12).Value = "Primera"
$ws.Cells.Item(2,13).Value = 120
$ws.Cells.Item(2,14).Value = 16000
$ws.Cells.Item(2,15).Value = 17000
$ws.Cells.Item(2,16).Value = 16500
$ws.Cells.Item(2,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(2,18).Value = "Región de O'Higgins"
$ws.Cells.Item(2,19).Value = 917
$ws.Cells.Item(2,20).Value = 18

# Row 3
$ws.Cells.Item(3,1).Value = 7
$ws.Cells.Item(3,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3,3).Value = "Ñuble"
$ws.Cells.Item(3,4).Value = 45096
$ws.Cells.Item(3,5).Value = 16
$ws.Cells.Item(3,6).Value = "Fruta"
$ws.Cells.Item(3,7).Value = 100107
$ws.Cells.Item(3,8).Value = "Otros"
$ws.Cells.Item(3,9).Value = 100107001
$ws.Cells.Item(3,10).Value = "Caqui"
$ws.Cells.Item(3,11).Value = "Mankaki"
$ws.Cells.Item(3,12).Value = "Primera"
$ws.Cells.Item(3,13).Value = 50
$ws.Cells.Item(3,14).Value = 18000
$ws.Cells.Item(3,15).Value = 18000
$ws.Cells.Item(3,16).Value = 18000
$ws.Cells.Item(3,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(3,18).Value = "Región del Maule"
$ws.Cells.Item(3,19).Value = 1000
$ws.Cells.Item(3,20).Value = 18

# Row 4
$ws.Cells.Item(4,1).Value = 7
$ws.Cells.Item(4,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4,3).Value = "Ñuble"
$ws.Cells.Item(4,4).Value = 44344
$ws.Cells.Item(4,5).Value = 16
$ws.Cells.Item(4,6).Value = "Fruta"
$ws.Cells.Item(4,7).Value = 100107
$ws.Cells.Item(4,8).Value = "Otros"
$ws.Cells.Item(4,9).Value = 100107001
$ws.Cells.Item(4,10).Value = "Caqui"
$ws.Cells.Item(4,11).Value = "Mankaki"
$ws.Cells.Item(4,12).Value = "Primera"
$ws.Cells.Item(4,13).Value = 120
$ws.Cells.Item(4,14).Value = 13000
$ws.Cells.Item(4,15).Value = 14000
$ws.Cells.Item(4,16).Value = 13500
$ws.Cells.Item(4,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(4,18).Value = "Provincia de Curicó"
$ws.Cells.Item(4,19).Value = 750
$ws.Cells.Item(4,20).Value = 18

# Row 5
$ws.Cells.Item(5,1).Value = 7
$ws.Cells.Item(5,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5,3).Value = "Ñuble"
$ws.Cells.Item(5,4).Value = 44330
$ws.Cells.Item(5,5).Value = 16
$ws.Cells.Item(5,6).Value = "Fruta"
$ws.Cells.Item(5,7).Value = 100107
$ws.Cells.Item(5,8).Value = "Otros"
$ws.Cells.Item(5,9).Value = 100107001
$ws.Cells.Item(5,10).Value = "Caqui"
$ws.Cells.Item(5,11).Value = "Mankaki"
$ws.Cells.Item(5,12).Value = "Primera"
$ws.Cells.Item(5,13).Value = 60
$ws.Cells.Item(5,14).Value = 15000
$ws.Cells.Item(5,15).Value = 16000
$ws.Cells.Item(5,16).Value = 15500
$ws.Cells.Item(5,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(5,18).Value = "Provincia de Curicó"
$ws.Cells.Item(5,19).Value = 861
$ws.Cells.Item(5,20).Value = 18

# Row 6
$ws.Cells.Item(6,1).Value = 7
$ws.Cells.Item(6,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6,3).Value = "Ñuble"
$ws.Cells.Item(6,4).Value = 44316
$ws.Cells.Item(6,5).Value = 16
$ws.Cells.Item(6,6).Value = "Fruta"
$ws.Cells.Item(6,7).Value = 100107
$ws.Cells.Item(6,8).Value = "Otros"
$ws.Cells.Item(6,9).Value = 100107001
$ws.Cells.Item(6,10).Value = "Caqui"
$ws.Cells.Item(6,11).Value = "Mankaki"
$ws.Cells.Item(6,12).Value = "Primera"
$ws.Cells.Item(6,13).Value = 60
$ws.Cells.Item(6,14).Value = 17500
$ws.Cells.Item(6,15).Value = 18000
$ws.Cells.Item(6,16).Value = 17750
$ws.Cells.Item(6,17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(6,18).Value = "Región de O'Higgins"
$ws.Cells.Item(6,19).Value = 1109
$ws.Cells.Item(6,20).Value = 16

# Row 7
$ws.Cells.Item(7,1).Value = 7
$ws.Cells.Item(7,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7,3).Value = "Ñuble"
$ws.Cells.Item(7,4).Value = 44316
$ws.Cells.Item(7,5).Value = 16
$ws.Cells.Item(7,6).Value = "Fruta"
$ws.Cells.Item(7,7).Value = 100107
$ws.Cells.Item(7,8).Value = "Otros"
$ws.Cells.Item(7,9).Value = 100107001
$ws.Cells.Item(7,10).Value = "Caqui"
$ws.Cells.Item(7,11).Value = "Mankaki"
$ws.Cells.Item(7,12).Value = "Segunda"
$ws.Cells.Item(7,13).Value = 40
$ws.Cells.Item(7,14).Value = 16000
$ws.Cells.Item(7,15).Value = 16000
$ws.Cells.Item(7,16).Value = 16000
$ws.Cells.Item(7,17).Value = "$/caja 16 kilos granel"
$ws.Cells.Item(7,18).Value = "Región de O'Higgins"
$ws.Cells.Item(7,19).Value = 1000
$ws.Cells.Item(7,20).Value = 16

# Row 8
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(8,3).Value = "Ñuble"
$ws.Cells.Item(8,4).Value = 45086
$ws.Cells.Item(8,5).Value = 16
$ws.Cells.Item(8,6).Value = "Fruta"
$ws.Cells.Item(8,7).Value = 100107
$ws.Cells.Item(8,8).Value = "Otros"
$ws.Cells.Item(8,9).Value = 100107001
$ws.Cells.Item(8,10).Value = "Caqui"
$ws.Cells.Item(8,11).Value = "Mankaki"
$ws.Cells.Item(8,12).Value = "Primera"
$ws.Cells.Item(8,13).Value = 30
$ws.Cells.Item(8,14).Value = 18000
$ws.Cells.Item(8,15).Value = 18000
$ws.Cells.Item(8,16).Value = 18000
$ws.Cells.Item(8,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(8,18).Value = "Región del Maule"
$ws.Cells.Item(8,19).Value = 1000
$ws.Cells.Item(8,20).Value = 18

# Row 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9,3).Value = "Ñuble"
$ws.Cells.Item(9,4).Value = 45092
$ws.Cells.Item(9,5).Value = 16
$ws.Cells.Item(9,6).Value = "Fruta"
$ws.Cells.Item(9,7).Value = 100107
$ws.Cells.Item(9,8).Value = "Otros"
$ws.Cells.Item(9,9).Value = 100107001
$ws.Cells.Item(9,10).Value = "Caqui"
$ws.Cells.Item(9,11).Value = "Mankaki"
$ws.Cells.Item(9,12).Value = "Primera"
$ws.Cells.Item(9,13).Value = 35
$ws.Cells.Item(9,14).Value = 18000
$ws.Cells.Item(9,15).Value = 19000
$ws.Cells.Item(9,16).Value = 18571
$ws.Cells.Item(9,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(9,18).Value = "Región del Maule"
$ws.Cells.Item(9,19).Value = 1032
$ws.Cells.Item(9,20).Value = 18

# Row 10
$ws.Cells.Item(10,1).Value = 7
$ws.Cells.Item(10,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10,3).Value = "Ñuble"
$ws.Cells.Item(10,4).Value = 45113
$ws.Cells.Item(10,5).Value = 16
$ws.Cells.Item(10,6).Value = "Fruta"
$ws.Cells.Item(10,7).Value = 100107
$ws.Cells.Item(10,8).Value = "Otros"
$ws.Cells.Item(10,9).Value = 100107001
$ws.Cells.Item(10,10).Value = "Caqui"
$ws.Cells.Item(10,11).Value = "Mankaki"
$ws.Cells.Item(10,12).Value = "Primera"
$ws.Cells.Item(10,13).Value = 20
$ws.Cells.Item(10,14).Value = 13000
$ws.Cells.Item(10,15).Value = 13000
$ws.Cells.Item(10,16).Value = 13000
$ws.Cells.Item(10,17).Value = "$/caja 12 kilos granel"
$ws.Cells.Item(10,18).Value = "Región del Maule"
$ws.Cells.Item(10,19).Value = 13000
$ws.Cells.Item(10,20).Value = 1

# Row 11
$ws.Cells.Item(11,1).Value = 7
$ws.Cells.Item(11,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(11,3).Value = "Ñuble"
$ws.Cells.Item(11,4).Value = 45107
$ws.Cells.Item(11,5).Value = 16
$ws.Cells.Item(11,6).Value = "Fruta"
$ws.Cells.Item(11,7).Value = 100107
$ws.Cells.Item(11,8).Value = "Otros"
$ws.Cells.Item(11,9).Value = 100107001
$ws.Cells.Item(11,10).Value = "Caqui"
$ws.Cells.Item(11,11).Value = "Mankaki"
$ws.Cells.Item(11,12).Value = "Primera"
$ws.Cells.Item(11,13).Value = 30
$ws.Cells.Item(11,14).Value = 18000
$ws.Cells.Item(11,15).Value = 18000
$ws.Cells.Item(11,16).Value = 18000
$ws.Cells.Item(11,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(11,18).Value = "Región del Maule"
$ws.Cells.Item(11,19).Value = 1000
$ws.Cells.Item(11,20).Value = 18

# Row 12
$ws.Cells.Item(12,1).Value = 7
$ws.Cells.Item(12,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12,3).Value = "Ñuble"
$ws.Cells.Item(12,4).Value = 44334
$ws.Cells.Item(12,5).Value = 16
$ws.Cells.Item(12,6).Value = "Fruta"
$ws.Cells.Item(12,7).Value = 100107
$ws.Cells.Item(12,8).Value = "Otros"
$ws.Cells.Item(12,9).Value = 100107001
$ws.Cells.Item(12,10).Value = "Caqui"
$ws.Cells.Item(12,11).Value = "Mankaki"
$ws.Cells.Item(12,12).Value = "Primera"
$ws.Cells.Item(12,13).Value = 120
$ws.Cells.Item(12,14).Value = 12000
$ws.Cells.Item(12,15).Value = 13000
$ws.Cells.Item(12,16).Value = 12500
$ws.Cells.Item(12,17).Value = "$/caja 12 kilos empedrada"
$ws.Cells.Item(12,18).Value = "Región de O'Higgins"
$ws.Cells.Item(12,19).Value = 1042
$ws.Cells.Item(12,20).Value = 12

# Row 13
$ws.Cells.Item(13,1).Value = 7
$ws.Cells.Item(13,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13,3).Value = "Ñuble"
$ws.Cells.Item(13,4).Value = 45085
$ws.Cells.Item(13,5).Value = 16
$ws.Cells.Item(13,6).Value = "Fruta"
$ws.Cells.Item(13,7).Value = 100107
$ws.Cells.Item(13,8).Value = "Otros"
$ws.Cells.Item(13,9).Value = 100107001
$ws.Cells.Item(13,10).Value = "Caqui"
$ws.Cells.Item(13,11).Value = "Mankaki"
$ws.Cells.Item(13,12).Value = "Primera"
$ws.Cells.Item(13,13).Value = 60
$ws.Cells.Item(13,14).Value = 18000
$ws.Cells.Item(13,15).Value = 19000
$ws.Cells.Item(13,16).Value = 18500
$ws.Cells.Item(13,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(13,18).Value = "Región del Maule"
$ws.Cells.Item(13,19).Value = 1028
$ws.Cells.Item(13,20).Value = 18

# Row 14
$ws.Cells.Item(14,1).Value = 7
$ws.Cells.Item(14,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14,3).Value = "Ñuble"
$ws.Cells.Item(14,4).Value = 45111
$ws.Cells.Item(14,5).Value = 16
$ws.Cells.Item(14,6).Value = "Fruta"
$ws.Cells.Item(14,7).Value = 100107
$ws.Cells.Item(14,8).Value = "Otros"
$ws.Cells.Item(14,9).Value = 100107001
$ws.Cells.Item(14,10).Value = "Caqui"
$ws.Cells.Item(14,11).Value = "Mankaki"
$ws.Cells.Item(14,12).Value = "Primera"
$ws.Cells.Item(14,13).Value = 30
$ws.Cells.Item(14,14).Value = 13000
$ws.Cells.Item(14,15).Value = 13000
$ws.Cells.Item(14,16).Value = 13000
$ws.Cells.Item(14,17).Value = "$/caja 12 kilos granel"
$ws.Cells.Item(14,18).Value = "Región del Maule"
$ws.Cells.Item(14,19).Value = 13000
$ws.Cells.Item(14,20).Value = 1

# Row 15
$ws.Cells.Item(15,1).Value = 7
$ws.Cells.Item(15,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15,3).Value = "Ñuble"
$ws.Cells.Item(15,4).Value = 45112
$ws.Cells.Item(15,5).Value = 16
$ws.Cells.Item(15,6).Value = "Fruta"
$ws.Cells.Item(15,7).Value = 100107
$ws.Cells.Item(15,8).Value = "Otros"
$ws.Cells.Item(15,9).Value = 100107001
$ws.Cells.Item(15,10).Value = "Caqui"
$ws.Cells.Item(15,11).Value = "Mankaki"
$ws.Cells.Item(15,12).Value = "Primera"
$ws.Cells.Item(15,13).Value = 30
$ws.Cells.Item(15,14).Value = 13000
$ws.Cells.Item(15,15).Value = 13000
$ws.Cells.Item(15,16).Value = 13000
$ws.Cells.Item(15,17).Value = "$/caja 12 kilos granel"
$ws.Cells.Item(15,18).Value = "Región del Maule"
$ws.Cells.Item(15,19).Value = 13000
$ws.Cells.Item(15,20).Value = 1

# Row 16
$ws.Cells.Item(16,1).Value = 7
$ws.Cells.Item(16,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16,3).Value = "Ñuble"
$ws.Cells.Item(16,4).Value = 45099
$ws.Cells.Item(16,5).Value = 16
$ws.Cells.Item(16,6).Value = "Fruta"
$ws.Cells.Item(16,7).Value = 100107
$ws.Cells.Item(16,8).Value = "Otros"
$ws.Cells.Item(16,9).Value = 100107001
$ws.Cells.Item(16,10).Value = "Caqui"
$ws.Cells.Item(16,11).Value = "Mankaki"
$ws.Cells.Item(16,12).Value = "Primera"
$ws.Cells.Item(16,13).Value = 50
$ws.Cells.Item(16,14).Value = 18000
$ws.Cells.Item(16,15).Value = 18000
$ws.Cells.Item(16,16).Value = 18000
$ws.Cells.Item(16,17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(16,18).Value = "Región del Maule"
$ws.Cells.Item(16,19).Value = 1000
$ws.Cells.Item(16,20).Value = 18

# Ensure the date column (D) keeps the date number format for the newly added rows
$ws.Cells.Item(14,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
